$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JAN-2021")
$ws.Activate()

# Insert three new rows (15:17) for the new Jan 25/27/29 log entries, pushing
# the legend/status rows (old 17-21) down to 20-24.
$ws.Rows("15:17").Insert()

# Copy formatting (styles) for columns A-F from the row above (row 14) onto
# the newly inserted rows, since a bare Insert() drops border formatting.
$ws.Range("A14:F14").Copy()
$ws.Range("A15:F17").PasteSpecial(-4122)

# Column G: rows 15 & 16 reuse the "wrap, centered, bordered" style used by
# column C (e.g. C2); row 17 reuses the G14 (wrap + border) style.
$ws.Range("C2").Copy()
$ws.Range("G15:G16").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("G17").PasteSpecial(-4122)

# --- Row 15: Jan 25 2021 ---
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 44221
$ws.Cells.Item(15, 3).Value = "Laptop SSD changing"
$ws.Cells.Item(15, 4).Value = "Laptop SSD changing"
$ws.Cells.Item(15, 7).Value = "Laptop Service"

# --- Row 16: Jan 27 2021 ---
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 44223
$ws.Cells.Item(16, 3).Value = "Laptop Software installed"
$ws.Cells.Item(16, 4).Value = "Laptop Software installed"
$ws.Cells.Item(16, 7).Value = "Laptop Software installed"

# --- Row 17: Jan 29 2021 ---
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 44225
$ws.Cells.Item(17, 3).Value = "QMVAR 2.0, Sony samsung testing"
$ws.Cells.Item(17, 7).Value = "Working View patr setup users"
$ws.Cells.Item(17, 4).Value = "Setup User View part"
$ws.Cells.Item(17, 5).Value = 0.8
$ws.Cells.Item(17, 6).Value = "Completed"

# Update the selection to reflect the author's last-saved cursor position.
$ws.Range("E17").Select()
